# "Output menu to spreadsheet."
# Adds a "winter" sub-block (columns I/J) to the two per-batch tables, adds a
# running total column (K) to both tables, and adds year/summer/winter
# average columns (L/M/N) with header labels + AVERAGE() formulas to both
# tables. Also extends the second ("actual batch") table with the C:I data
# that was missing, plus matching SUM()/formula columns, and moves the
# viewport down a bit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New column widths for the L/M/N (year/summer/winter average) columns.
# ---------------------------------------------------------------------
$ws.Columns.Item(12).ColumnWidth = 11.77734375
$ws.Columns.Item(13).ColumnWidth = 11.6640625
$ws.Columns.Item(14).ColumnWidth = 10

# ---------------------------------------------------------------------
# First table (rows 161-170): "winter" header row + running total (K) +
# year/summer/winter average block (L/M/N).
# ---------------------------------------------------------------------

# Row 162 ("winter" sub-header over I:J, the last two -- winter -- months)
$ws.Cells.Item(164, 11).Value = "total for 370 days"
$ws.Cells.Item(162, 9).Value = "winter"
$ws.Cells.Item(162, 10).Value = "winter"

# Running total column K for rows 165-167
$ws.Cells.Item(165, 11).Formula = "=SUM(B165:J165)"
$ws.Cells.Item(166, 11).Formula = "=SUM(B166:J166)"
$ws.Cells.Item(167, 11).Formula = "=SUM(B167:J167)"

# Header labels L167:N167
$ws.Cells.Item(167, 12).Value = "year average"
$ws.Cells.Item(167, 14).Value = "winter average"
$ws.Cells.Item(167, 13).Value = "summer average"
$ws.Cells.Item(167, 12).HorizontalAlignment = -4108

# Average formulas L168:N170
$ws.Cells.Item(168, 12).Formula = "=AVERAGE(B168:J168)"
$ws.Cells.Item(168, 13).Formula = "=AVERAGE(B168:H168)"
$ws.Cells.Item(168, 14).Formula = "=AVERAGE(I168:J168)"

$ws.Cells.Item(169, 12).Formula = "=AVERAGE(B169:J169)"
$ws.Cells.Item(169, 13).Formula = "=AVERAGE(B169:H169)"
$ws.Cells.Item(169, 14).Formula = "=AVERAGE(I169:J169)"

$ws.Cells.Item(170, 12).Formula = "=AVERAGE(B170:J170)"
$ws.Cells.Item(170, 13).Formula = "=AVERAGE(B170:H170)"
$ws.Cells.Item(170, 14).Formula = "=AVERAGE(I170:J170)"

# ---------------------------------------------------------------------
# Second table (rows 172-180): same shape, plus it previously only had
# columns B and J filled in -- fill in the missing C:I months too.
# ---------------------------------------------------------------------

# Row 172 "winter" sub-header over I:J
$ws.Cells.Item(172, 9).Value = "winter"
$ws.Cells.Item(172, 10).Value = "winter"

# Row 173 (batch numbers) C:I
$ws.Cells.Item(173, 3).Value = 2
$ws.Cells.Item(173, 4).Value = 3
$ws.Cells.Item(173, 5).Value = 4
$ws.Cells.Item(173, 6).Value = 5
$ws.Cells.Item(173, 7).Value = 6
$ws.Cells.Item(173, 8).Value = 7
$ws.Cells.Item(173, 9).Value = 8

# Row 174 (Matrix size) C:I + running-total header K174
$ws.Cells.Item(174, 3).Value = 3149
$ws.Cells.Item(174, 4).Value = 3054
$ws.Cells.Item(174, 5).Value = 3093
$ws.Cells.Item(174, 6).Value = 3154
$ws.Cells.Item(174, 7).Value = 3077
$ws.Cells.Item(174, 8).Value = 3026
$ws.Cells.Item(174, 9).Value = 3012
$ws.Cells.Item(174, 11).Value = "total for 370 days"

# Row 175 (cost £) C,D,E,F,(no G),H,I + K175 total
$ws.Cells.Item(175, 3).Value = 28349
$ws.Cells.Item(175, 4).Value = 27874
$ws.Cells.Item(175, 5).Value = 30408
$ws.Cells.Item(175, 6).Value = 30552
$ws.Cells.Item(175, 8).Value = 27287
$ws.Cells.Item(175, 9).Value = 29005
$ws.Cells.Item(175, 11).Formula = "=SUM(B175:J175)"

# Row 176 (emissions kg) C,D,E,F,(no G),H,I + K176 total
$ws.Cells.Item(176, 3).Value = 23656
$ws.Cells.Item(176, 4).Value = 24225
$ws.Cells.Item(176, 5).Value = 26364
$ws.Cells.Item(176, 6).Value = 26706
$ws.Cells.Item(176, 8).Value = 24405
$ws.Cells.Item(176, 9).Value = 22464
$ws.Cells.Item(176, 11).Formula = "=SUM(B176:J176)"

# Row 177 (food waste sum nutrients) C,D,E,F,(no G),H,I + K177 total + headers
$ws.Cells.Item(177, 3).Value = 836547
$ws.Cells.Item(177, 4).Value = 844315
$ws.Cells.Item(177, 5).Value = 800038
$ws.Cells.Item(177, 6).Value = 968130
$ws.Cells.Item(177, 8).Value = 1215430
$ws.Cells.Item(177, 9).Value = 1061965
$ws.Cells.Item(177, 11).Formula = "=SUM(B177:J177)"

$ws.Cells.Item(177, 12).Value = "year average"
$ws.Cells.Item(177, 14).Value = "winter average"
$ws.Cells.Item(177, 13).Value = "summer average"
$ws.Cells.Item(177, 12).HorizontalAlignment = -4108

# Row 178 (cost/(days x people)) C,D,E,F,(no G),H,I + L/M/N averages
$ws.Cells.Item(178, 3).Formula = "=28349/3149"
$ws.Cells.Item(178, 4).Formula = "=27874/3054"
$ws.Cells.Item(178, 5).Formula = "=30408/3093"
$ws.Cells.Item(178, 6).Formula = "=30552/3154"
$ws.Cells.Item(178, 8).Formula = "=27287/3026"
$ws.Cells.Item(178, 9).Formula = "=29005/3012"
$ws.Cells.Item(178, 12).Formula = "=AVERAGE(B178:J178)"
$ws.Cells.Item(178, 13).Formula = "=AVERAGE(B178:H178)"
$ws.Cells.Item(178, 14).Formula = "=AVERAGE(I178:J178)"

# Row 179 (emissions /(days x people)) C,D,E,F,(no G),H,I + L/M/N averages
$ws.Cells.Item(179, 3).Formula = "=23656/3149"
$ws.Cells.Item(179, 4).Formula = "=24225/3054"
$ws.Cells.Item(179, 5).Formula = "=26364/3093"
$ws.Cells.Item(179, 6).Formula = "=26706/3154"
$ws.Cells.Item(179, 8).Formula = "=24405/3026"
$ws.Cells.Item(179, 9).Formula = "=22464/3012"
$ws.Cells.Item(179, 12).Formula = "=AVERAGE(B179:J179)"
$ws.Cells.Item(179, 13).Formula = "=AVERAGE(B179:H179)"
$ws.Cells.Item(179, 14).Formula = "=AVERAGE(I179:J179)"

# Row 180 (food waste /(days x people)) C,D,E,F,(no G),H,I + L/M/N averages
$ws.Cells.Item(180, 3).Formula = "=836547/3149"
$ws.Cells.Item(180, 4).Formula = "=844315/3054"
$ws.Cells.Item(180, 5).Formula = "=800038/3093"
$ws.Cells.Item(180, 6).Formula = "=968130/3154"
$ws.Cells.Item(180, 8).Formula = "=1215430/3026"
$ws.Cells.Item(180, 9).Formula = "=1061965/3012"
$ws.Cells.Item(180, 12).Formula = "=AVERAGE(B180:J180)"
$ws.Cells.Item(180, 13).Formula = "=AVERAGE(B180:H180)"
$ws.Cells.Item(180, 14).Formula = "=AVERAGE(I180:J180)"

# ---------------------------------------------------------------------
# Move the viewport: scroll down a bit and move the active selection to
# O190 (just past the new data), matching where the author ended up.
# ---------------------------------------------------------------------
$excel.Goto($ws.Range("A165"), $true)
$ws.Range("O190").Select()
